# Remove "Table 4 - Research Parameters of Study" and
# "Table 5 - Outcome Measurements of Study" sections (heading, spacer
# paragraph, the data table itself, and the trailing page break) from
# the document, along with the reviewer comment that lived inside
# Table 5.

$d = $word.ActiveDocument

function Find-ParagraphIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -match [regex]::Escape($text)) {
            return $i
        }
    }
    return -1
}

function Find-TableIndexContainingCellText($text) {
    for ($i = 1; $i -le $d.Tables.Count; $i++) {
        if ($d.Tables($i).Range.Text -match [regex]::Escape($text)) {
            return $i
        }
    }
    return -1
}

# Locate the two tables to remove by distinctive header-cell text
# ("Design" only appears in Table 4, "Measurement(s)" only in Table 5).
$table4Idx = Find-TableIndexContainingCellText("Design")
$table5Idx = Find-TableIndexContainingCellText("Measurement(s)")

# Delete the higher-indexed table first so the other index stays valid.
if ($table5Idx -gt $table4Idx) {
    $d.Tables($table5Idx).Delete()
    $d.Tables($table4Idx).Delete()
} else {
    $d.Tables($table4Idx).Delete()
    $d.Tables($table5Idx).Delete()
}

# Remove the now-empty heading / spacer / page-break paragraphs that
# used to surround those two tables: everything from the "Table 4"
# heading paragraph up to (but not including) the "Table 6" heading
# paragraph.
$startIdx = Find-ParagraphIndexByText("Table 4")
$endHeadingIdx = Find-ParagraphIndexByText("Table 6")
$endIdx = $endHeadingIdx - 1

$pStart = $d.Paragraphs($startIdx)
$pEnd = $d.Paragraphs($endIdx)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.Delete()

# Drop the reviewer comment that was anchored inside the deleted
# Table 5 ("consider describing this to students").
$commentCount = $d.Comments.Count
for ($i = $commentCount; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}

Write-Output "Tables remaining: $($d.Tables.Count)"
Write-Output "Comments remaining: $($d.Comments.Count)"
